$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 0.1887967460050058
$ws.Range("C2").Value = 2.743091716767341
$ws.Range("D2").Value = 23.88233659257806
$ws.Range("E2").Value = 4.886955759220464
$ws.Range("F2").Value = 4.931898846174919
$ws.Range("G2").Value = 51

# Row 3
$ws.Range("B3").Value = 0.352110747706186
$ws.Range("C3").Value = 2.607870699580363
$ws.Range("D3").Value = 22.56754235743245
$ws.Range("E3").Value = 4.750530744814989
$ws.Range("F3").Value = 4.785560805511427
$ws.Range("G3").Value = 50

# Row 4
$ws.Range("B4").Value = 0.2174705470141381
$ws.Range("C4").Value = 2.756943539436507
$ws.Range("D4").Value = 23.51947414538977
$ws.Range("E4").Value = 4.849688046193257
$ws.Range("F4").Value = 4.895016289175286
$ws.Range("G4").Value = 49

# Row 5
$ws.Range("B5").Value = 0.4431031871169503
$ws.Range("C5").Value = 2.719424870216065
$ws.Range("D5").Value = 23.17433994440393
$ws.Range("E5").Value = 4.813973405037041
$ws.Range("F5").Value = 4.844263939605876
$ws.Range("G5").Value = 48

# Row 6
$ws.Range("B6").Value = 0.2052815252275707
$ws.Range("C6").Value = 2.707768026847837
$ws.Range("D6").Value = 24.13652786648203
$ws.Range("E6").Value = 4.912894041853745
$ws.Range("F6").Value = 4.961670927367957
$ws.Range("G6").Value = 47

# Row 7
$ws.Range("B7").Value = 0.4590069508388199
$ws.Range("C7").Value = 2.706948687152313
$ws.Range("D7").Value = 24.1233625610406
$ws.Range("E7").Value = 4.911553986371381
$ws.Range("F7").Value = 4.944094250912167
$ws.Range("G7").Value = 46

# Row 8
$ws.Range("B8").Value = 0.04636964274274196
$ws.Range("C8").Value = 2.422883991970613
$ws.Range("D8").Value = 21.4973978904842
$ws.Range("E8").Value = 4.636528646572153
$ws.Range("F8").Value = 4.688685967794817
$ws.Range("G8").Value = 45

# Row 9
$ws.Range("B9").Value = 0.2175654990191061
$ws.Range("C9").Value = 2.393845526388351
$ws.Range("D9").Value = 21.18684633999888
$ws.Range("E9").Value = 4.60291715545684
$ws.Range("F9").Value = 4.650927664705682
$ws.Range("G9").Value = 44

# Row 10
$ws.Range("B10").Value = 0.02747142410952863
$ws.Range("C10").Value = 2.403862371068361
$ws.Range("D10").Value = 21.54546590760865
$ws.Range("E10").Value = 4.641709373453777
$ws.Range("F10").Value = 4.696560501412658
$ws.Range("G10").Value = 43

# Row 11
$ws.Range("B11").Value = 0.2558878751990361
$ws.Range("C11").Value = 2.356132134174123
$ws.Range("D11").Value = 21.26811481140563
$ws.Range("E11").Value = 4.611736637255605
$ws.Range("F11").Value = 4.660447797710927
$ws.Range("G11").Value = 42
